$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "PAY31F5C6BB"
$ws.Range("A3").Value = "PAY418F7842"
$ws.Range("A4").Value = "PAYF64E519F"
$ws.Range("A5").Value = "PAYC0805C7B"
$ws.Range("A6").Value = "PAY31C994D7"
$ws.Range("A7").Value = "PAYA271F21A"
$ws.Range("A8").Value = "PAY65C2DA42"
$ws.Range("A9").Value = "PAYF147FED1"
$ws.Range("A10").Value = "PAY3EDC62C4"
$ws.Range("A11").Value = "PAY0753339B"
$ws.Range("A12").Value = "PAYC2988104"
$ws.Range("A13").Value = "PAY058E3389"
$ws.Range("A14").Value = "PAYC14714F7"
$ws.Range("A15").Value = "PAY9AD4C4C7"
$ws.Range("A16").Value = "PAYEB359F50"
$ws.Range("A17").Value = "PAY18113F7D"
$ws.Range("A18").Value = "PAYA9360529"
$ws.Range("A19").Value = "PAY0AD7D2FA"
$ws.Range("A20").Value = "PAYDF38A650"
$ws.Range("A21").Value = "PAYDC17523D"
$ws.Range("A22").Value = "PAY52CECD71"
$ws.Range("A23").Value = "PAYDEBD214A"
$ws.Range("A24").Value = "PAY39F8FB60"
$ws.Range("A25").Value = "PAY7DAEB3BE"
$ws.Range("A26").Value = "PAYDD8BB0DD"
$ws.Range("A27").Value = "PAY0A399CC5"
$ws.Range("A28").Value = "PAY548E68E8"
$ws.Range("A29").Value = "PAY7C16C19E"
$ws.Range("A30").Value = "PAYD2C29BB6"
$ws.Range("A31").Value = "PAY859E9C9F"
$ws.Range("A32").Value = "PAY8B6EFD19"
$ws.Range("A33").Value = "PAY00536CDD"
$ws.Range("A34").Value = "PAY9CED8223"
$ws.Range("A35").Value = "PAY130384CE"
$ws.Range("A36").Value = "PAYBDFD264C"
$ws.Range("A37").Value = "PAYA1AD1F96"
$ws.Range("A38").Value = "PAY804D0535"
$ws.Range("A39").Value = "PAY4376EE21"
$ws.Range("A40").Value = "PAYFFE12DEE"
$ws.Range("A41").Value = "PAYC48E5638"
$ws.Range("A42").Value = "PAYFA15388F"
$ws.Range("A43").Value = "PAY9E313304"
$ws.Range("A44").Value = "PAY5A09C97B"
$ws.Range("A45").Value = "PAYD992370E"
$ws.Range("A46").Value = "PAYDCC1B970"
$ws.Range("A47").Value = "PAYB2ECAE66"
$ws.Range("A48").Value = "PAY94352C8F"
$ws.Range("A49").Value = "PAYE2EC6854"
$ws.Range("A50").Value = "PAY051EF33F"
$ws.Range("A51").Value = "PAYFF6FFC98"
$ws.Range("A52").Value = "PAYB137875B"
$ws.Range("A53").Value = "PAY54198F3A"
$ws.Range("A54").Value = "PAY43B29987"
$ws.Range("A55").Value = "PAYABB40C41"
$ws.Range("A56").Value = "PAY0F9FC53D"
$ws.Range("A57").Value = "PAY3B6B1008"
$ws.Range("A58").Value = "PAY6ABA67FC"
$ws.Range("A59").Value = "PAYB7DF0AD9"
$ws.Range("A60").Value = "PAYD46410D0"
$ws.Range("A61").Value = "PAYE0B451D4"
$ws.Range("A62").Value = "PAY39CD28D1"
$ws.Range("A63").Value = "PAY17B6D53B"
$ws.Range("A64").Value = "PAY91D4EC6C"
$ws.Range("A65").Value = "PAY2E74B4C0"
$ws.Range("A66").Value = "PAY489A0D5B"
$ws.Range("A67").Value = "PAY8639090C"
$ws.Range("A68").Value = "PAYF9BC1405"
$ws.Range("A69").Value = "PAYEA6C7013"
$ws.Range("A70").Value = "PAY8A5D611C"
$ws.Range("A71").Value = "PAYBB59EC8C"
$ws.Range("A72").Value = "PAYF092C6DF"
$ws.Range("A73").Value = "PAY7A2E1571"
$ws.Range("A74").Value = "PAY245318BD"
$ws.Range("A75").Value = "PAY672B91C2"
$ws.Range("A76").Value = "PAYE32022FB"
$ws.Range("A77").Value = "PAYCA329B7F"
$ws.Range("A78").Value = "PAY1121D59A"
$ws.Range("A79").Value = "PAYD6CC8647"
$ws.Range("A80").Value = "PAY3C641919"
$ws.Range("A81").Value = "PAY3375B291"
$ws.Range("A82").Value = "PAY16E69662"
$ws.Range("A83").Value = "PAY2278B2B3"
$ws.Range("A84").Value = "PAYC414902D"
$ws.Range("A85").Value = "PAY46EBF25E"
$ws.Range("A86").Value = "PAY014B512A"
$ws.Range("A87").Value = "PAY9A9F5D6D"
$ws.Range("A88").Value = "PAY9B005484"
$ws.Range("A89").Value = "PAY1E6796F4"
$ws.Range("A90").Value = "PAYC93E0F26"
$ws.Range("A91").Value = "PAY1BFB5C09"
$ws.Range("A92").Value = "PAY4BFC7128"
$ws.Range("A93").Value = "PAY7B2D146C"
$ws.Range("A94").Value = "PAY6274B391"
$ws.Range("A95").Value = "PAYE0DC9DAA"
$ws.Range("A96").Value = "PAYC8D89FDC"
$ws.Range("A97").Value = "PAY547586C9"
$ws.Range("A98").Value = "PAYEAA1A7EA"
$ws.Range("A99").Value = "PAY9BD5E633"
$ws.Range("A100").Value = "PAY646FF52A"
$ws.Range("A101").Value = "PAYB7F59EC0"
$ws.Range("A102").Value = "PAY888AFF54"
$ws.Range("A103").Value = "PAY9B7727B4"
$ws.Range("A104").Value = "PAYFB32B569"
$ws.Range("A105").Value = "PAY600FAC7C"
$ws.Range("A106").Value = "PAY0983A179"
$ws.Range("A107").Value = "PAY36F32F07"
$ws.Range("A108").Value = "PAYF6CA6A3E"
$ws.Range("A109").Value = "PAYF2975F35"
$ws.Range("A110").Value = "PAY6933A2AF"
$ws.Range("A111").Value = "PAY35331E86"
$ws.Range("A112").Value = "PAYD13FF7F9"
$ws.Range("A113").Value = "PAYFC06E3EA"
$ws.Range("A114").Value = "PAY1EACECD2"
$ws.Range("A115").Value = "PAY644C3631"
$ws.Range("A116").Value = "PAY1231AEC5"
$ws.Range("A117").Value = "PAYD7853B67"
$ws.Range("A118").Value = "PAY333B542B"
$ws.Range("A119").Value = "PAYC539950A"
$ws.Range("A120").Value = "PAY77788F76"
$ws.Range("A121").Value = "PAY063FF8F1"
$ws.Range("A122").Value = "PAY44DCA5C4"
$ws.Range("A123").Value = "PAY3ED174BA"
$ws.Range("A124").Value = "PAY04913EAB"
$ws.Range("A125").Value = "PAY52762535"
$ws.Range("A126").Value = "PAY54C599F5"
$ws.Range("A127").Value = "PAYA7669A90"
$ws.Range("A128").Value = "PAYA2233404"
$ws.Range("A129").Value = "PAY9C1B0390"
$ws.Range("A130").Value = "PAY53286B38"
$ws.Range("A131").Value = "PAY555BB008"
$ws.Range("A132").Value = "PAY26C074D7"
$ws.Range("A133").Value = "PAYEF2A6B1B"
$ws.Range("A134").Value = "PAY0208B75F"
$ws.Range("A135").Value = "PAYE1E41FC5"
$ws.Range("A136").Value = "PAY1E7F1206"
$ws.Range("A137").Value = "PAY4EBEA434"
$ws.Range("A138").Value = "PAYC346B081"
$ws.Range("A139").Value = "PAY1AD6286A"
$ws.Range("A140").Value = "PAY2DC0F188"
$ws.Range("A141").Value = "PAY9ED810E3"
$ws.Range("A142").Value = "PAY3FF21238"
$ws.Range("A143").Value = "PAY1027E23A"
$ws.Range("A144").Value = "PAY8066B0DD"
$ws.Range("A145").Value = "PAYBB7D7098"
$ws.Range("A146").Value = "PAY18A5CB28"
$ws.Range("A147").Value = "PAY983A7ED9"
$ws.Range("A148").Value = "PAYB9F97E76"
$ws.Range("A149").Value = "PAY3E4445B8"
$ws.Range("A150").Value = "PAY087EDDDB"
$ws.Range("A151").Value = "PAY9511E35B"
$ws.Range("A152").Value = "PAYE3A10111"
$ws.Range("A153").Value = "PAY71612CBC"
$ws.Range("A154").Value = "PAY0B3628E9"
$ws.Range("A155").Value = "PAY64654A1B"
$ws.Range("A156").Value = "PAY84C7EE32"
$ws.Range("A157").Value = "PAYC2C0FD8B"
$ws.Range("A158").Value = "PAYA2B652AA"
$ws.Range("A159").Value = "PAY0BB76043"
$ws.Range("A160").Value = "PAY42660D67"
$ws.Range("A161").Value = "PAY68ED1447"
$ws.Range("A162").Value = "PAYB084C8B0"
$ws.Range("A163").Value = "PAYA448A10D"
$ws.Range("A164").Value = "PAYC6B1D147"
$ws.Range("A165").Value = "PAY8A050058"
$ws.Range("A166").Value = "PAY7B5E7F0A"
$ws.Range("A167").Value = "PAY6809A348"
$ws.Range("A168").Value = "PAY3994B353"
$ws.Range("A169").Value = "PAY3CEC20E3"
$ws.Range("A170").Value = "PAY754BB3F0"
$ws.Range("A171").Value = "PAYB3B38766"
$ws.Range("A172").Value = "PAYA48173CB"
$ws.Range("A173").Value = "PAY3BEEAA44"
$ws.Range("A174").Value = "PAY14E98E83"
$ws.Range("A175").Value = "PAYC8D21752"
$ws.Range("A176").Value = "PAY143B9A68"
$ws.Range("A177").Value = "PAY9863B98E"
$ws.Range("A178").Value = "PAY4F8750F6"
$ws.Range("A179").Value = "PAY17640301"
$ws.Range("A180").Value = "PAYC6C04594"
$ws.Range("A181").Value = "PAYDDC96F61"
$ws.Range("A182").Value = "PAYFADBD5B3"
$ws.Range("A183").Value = "PAY2A28EDE1"
$ws.Range("A184").Value = "PAYA22E9CB4"
$ws.Range("A185").Value = "PAY6FF055D9"
$ws.Range("A186").Value = "PAYFBC523E7"
$ws.Range("A187").Value = "PAY3A6624B6"
$ws.Range("A188").Value = "PAYE4EAA263"
$ws.Range("A189").Value = "PAYCA127CC7"
$ws.Range("A190").Value = "PAY4AB54403"
$ws.Range("A191").Value = "PAY80D02996"
$ws.Range("A192").Value = "PAY9EB84E82"
$ws.Range("A193").Value = "PAY7EF15669"
$ws.Range("A194").Value = "PAYB63E5444"
$ws.Range("A195").Value = "PAY2C969615"
$ws.Range("A196").Value = "PAY01634F2F"
$ws.Range("A197").Value = "PAY1CDE164E"
$ws.Range("A198").Value = "PAY39C7122C"
$ws.Range("A199").Value = "PAY8D1C563C"
$ws.Range("A200").Value = "PAYADA61C46"
$ws.Range("A201").Value = "PAY4B55DE3E"
$ws.Range("A202").Value = "PAYB03E12C6"
$ws.Range("A203").Value = "PAY4E18F03B"
$ws.Range("A204").Value = "PAY91489507"
$ws.Range("A205").Value = "PAY1F6DF16E"
$ws.Range("A206").Value = "PAYB095CD63"
$ws.Range("A207").Value = "PAYB5D26FAC"
$ws.Range("A208").Value = "PAY8F63FBFE"
$ws.Range("A209").Value = "PAY193304E5"
$ws.Range("A210").Value = "PAY7948C424"
$ws.Range("A211").Value = "PAYE25E8924"
$ws.Range("A212").Value = "PAY45864395"
$ws.Range("A213").Value = "PAY118E2F59"
$ws.Range("A214").Value = "PAY59B4EDF8"
$ws.Range("A215").Value = "PAY15EE2FBA"
$ws.Range("A216").Value = "PAY0B3F25DB"
$ws.Range("A217").Value = "PAYFC021E7B"
$ws.Range("A218").Value = "PAYA8453BBF"
$ws.Range("A219").Value = "PAY66DA74BE"
$ws.Range("A220").Value = "PAY67BD6706"
$ws.Range("A221").Value = "PAY6D61B60C"
$ws.Range("A222").Value = "PAYD3653AA3"
$ws.Range("A223").Value = "PAY97CCD90B"
$ws.Range("A224").Value = "PAYF195611F"
$ws.Range("A225").Value = "PAY158E39A1"
$ws.Range("A226").Value = "PAY10F7D20A"
$ws.Range("A227").Value = "PAY11D91F70"
$ws.Range("A228").Value = "PAY23C2495E"
$ws.Range("A229").Value = "PAYFF5A939E"
$ws.Range("A230").Value = "PAY8E4B7712"
$ws.Range("A231").Value = "PAY78BBD104"
$ws.Range("A232").Value = "PAY9722676D"
$ws.Range("A233").Value = "PAYB02AE1E6"
$ws.Range("A234").Value = "PAYA0C0D2DC"
$ws.Range("A235").Value = "PAY0806F7DC"
$ws.Range("A236").Value = "PAY1993B614"
$ws.Range("A237").Value = "PAY575686E9"
$ws.Range("A238").Value = "PAYFB85A403"
$ws.Range("A239").Value = "PAYAB270367"
$ws.Range("A240").Value = "PAYDDA62D20"
$ws.Range("A241").Value = "PAY87AA226F"
$ws.Range("A242").Value = "PAY2AC5F3B5"
$ws.Range("A243").Value = "PAY85ACB2A8"
$ws.Range("A244").Value = "PAYC7202EF9"
$ws.Range("A245").Value = "PAY68E83CBF"
$ws.Range("A246").Value = "PAY40522A08"
$ws.Range("A247").Value = "PAY5D7967C9"
$ws.Range("A248").Value = "PAY92FCB4E2"
$ws.Range("A249").Value = "PAY86AE0F0B"
$ws.Range("A250").Value = "PAY3B5DB729"
$ws.Range("A251").Value = "PAYC6F15225"
$ws.Range("A252").Value = "PAY27BC5734"
$ws.Range("A253").Value = "PAY2BB8BD8E"
$ws.Range("A254").Value = "PAYF12C368C"
$ws.Range("A255").Value = "PAY410007DA"
$ws.Range("A256").Value = "PAY80042F19"
$ws.Range("A257").Value = "PAY4C1D4528"
$ws.Range("A258").Value = "PAYABD5064B"
$ws.Range("A259").Value = "PAY6ACF1833"
$ws.Range("A260").Value = "PAY7A479403"
$ws.Range("A261").Value = "PAYBD71FCC1"
$ws.Range("A262").Value = "PAYB7BBE7EE"
$ws.Range("A263").Value = "PAY4C933AF7"
$ws.Range("A264").Value = "PAYC2232362"
$ws.Range("A265").Value = "PAY91B7F2E9"
$ws.Range("A266").Value = "PAY1416A0D2"
$ws.Range("A267").Value = "PAY2D969A6E"
$ws.Range("A268").Value = "PAYE9D0763E"
$ws.Range("A269").Value = "PAY254D58C3"
$ws.Range("A270").Value = "PAY1F6A6CDC"
$ws.Range("A271").Value = "PAY06FE5420"
$ws.Range("A272").Value = "PAY72C5B855"
$ws.Range("A273").Value = "PAY3678E941"
$ws.Range("A274").Value = "PAY6AB82F0A"
$ws.Range("A275").Value = "PAYA20423D5"
$ws.Range("A276").Value = "PAYB831A4B3"
$ws.Range("A277").Value = "PAY6C08BAE9"
$ws.Range("A278").Value = "PAY4B8151A0"
$ws.Range("A279").Value = "PAYE2C87301"
$ws.Range("A280").Value = "PAY2DAECF46"
$ws.Range("A281").Value = "PAY8FB6EC12"
$ws.Range("A282").Value = "PAY30A0707B"
$ws.Range("A283").Value = "PAYA8A56D0B"
$ws.Range("A284").Value = "PAY721EAFE7"
$ws.Range("A285").Value = "PAYE723AA85"
$ws.Range("A286").Value = "PAY6788E4E4"
$ws.Range("A287").Value = "PAY4EE4C9D0"
$ws.Range("A288").Value = "PAY008CC26A"
$ws.Range("A289").Value = "PAYFC9B1E2F"
$ws.Range("A290").Value = "PAY669B79D1"
$ws.Range("A291").Value = "PAYB3C231D7"
$ws.Range("A292").Value = "PAY0CC9A67A"
$ws.Range("A293").Value = "PAY15923326"
$ws.Range("A294").Value = "PAY31A539C5"
$ws.Range("A295").Value = "PAY63F10960"
$ws.Range("A296").Value = "PAYEF0EAAA0"
$ws.Range("A297").Value = "PAY3AE50DC8"
$ws.Range("A298").Value = "PAY2794B31D"
$ws.Range("A299").Value = "PAY8D5FF5AF"
$ws.Range("A300").Value = "PAY322DAB57"
$ws.Range("A301").Value = "PAY59C8007C"
$ws.Range("A302").Value = "PAY906961E1"
$ws.Range("A303").Value = "PAYBC7A7EB3"
$ws.Range("A304").Value = "PAY0E8C22F9"
$ws.Range("A305").Value = "PAY980A9BE4"
$ws.Range("A306").Value = "PAY53B6F7B2"
$ws.Range("A307").Value = "PAY6062369F"
$ws.Range("A308").Value = "PAY8F78A1E1"
$ws.Range("A309").Value = "PAYDA8A40A9"
$ws.Range("A310").Value = "PAYE1245297"
$ws.Range("A311").Value = "PAYA830E34C"
$ws.Range("A312").Value = "PAY526865A7"
$ws.Range("A313").Value = "PAY18C76C5B"
$ws.Range("A314").Value = "PAYFA25E0B9"
$ws.Range("A315").Value = "PAY4B20B015"
$ws.Range("A316").Value = "PAYD74F5E95"
$ws.Range("A317").Value = "PAY8EFCE3E4"
$ws.Range("A318").Value = "PAY93C334E5"
$ws.Range("A319").Value = "PAY5C8D7FFE"
$ws.Range("A320").Value = "PAYBC84B5E5"
$ws.Range("A321").Value = "PAY238A37A4"
$ws.Range("A322").Value = "PAY753D3B1F"
$ws.Range("A323").Value = "PAY1AB23028"
$ws.Range("A324").Value = "PAYB53B9C70"
$ws.Range("A325").Value = "PAY97850A24"
$ws.Range("A326").Value = "PAY4BF268D0"
$ws.Range("A327").Value = "PAYA73B0E97"
$ws.Range("A328").Value = "PAY9898838A"
$ws.Range("A329").Value = "PAYAB147055"
$ws.Range("A330").Value = "PAY28A1B19A"
$ws.Range("A331").Value = "PAYFABA92B1"
$ws.Range("A332").Value = "PAY879D1BF5"
$ws.Range("A333").Value = "PAY8ADB5E20"
$ws.Range("A334").Value = "PAY44E811A0"
$ws.Range("A335").Value = "PAY6EF784FD"
$ws.Range("A336").Value = "PAY09EF3E02"
$ws.Range("A337").Value = "PAY25AB5579"
$ws.Range("A338").Value = "PAY563DF2CE"
$ws.Range("A339").Value = "PAY85AFBDF1"
$ws.Range("A340").Value = "PAY6A624FF0"
$ws.Range("A341").Value = "PAYEE44C4BA"
$ws.Range("A342").Value = "PAY31C19612"
$ws.Range("A343").Value = "PAY5364648C"
$ws.Range("A344").Value = "PAY191F9B94"
$ws.Range("A345").Value = "PAY3C47FF71"
$ws.Range("A346").Value = "PAY2F2258C9"
$ws.Range("A347").Value = "PAYA1CC80BB"
$ws.Range("A348").Value = "PAYCB5E7228"
$ws.Range("A349").Value = "PAY6BAC2E4A"
$ws.Range("A350").Value = "PAY26F8900C"
$ws.Range("A351").Value = "PAY0A8255EA"
$ws.Range("A352").Value = "PAYB4BD7C1B"
$ws.Range("A353").Value = "PAY4A4DC978"
$ws.Range("A354").Value = "PAYA6840BA1"
$ws.Range("A355").Value = "PAY9A03ACB5"
$ws.Range("A356").Value = "PAYF7C87017"
$ws.Range("A357").Value = "PAY981C5985"
$ws.Range("A358").Value = "PAY62634E7D"
$ws.Range("A359").Value = "PAYB57E3C6F"
$ws.Range("A360").Value = "PAYC7DA8FE7"
$ws.Range("A361").Value = "PAY424C53A1"
$ws.Range("A362").Value = "PAYFBDED9D8"
$ws.Range("A363").Value = "PAYC2252E94"
$ws.Range("A364").Value = "PAYA6881D1E"
$ws.Range("A365").Value = "PAYA61D4217"
$ws.Range("A366").Value = "PAY4188FDCA"
$ws.Range("A367").Value = "PAYD1284A53"
$ws.Range("A368").Value = "PAYAF40AAF6"
$ws.Range("A369").Value = "PAYE8DB2122"
$ws.Range("A370").Value = "PAY5E56AA70"
$ws.Range("A371").Value = "PAY2F2D72E2"
$ws.Range("A372").Value = "PAY891BC671"
$ws.Range("A373").Value = "PAY604905AD"
$ws.Range("A374").Value = "PAY7DFFA931"
$ws.Range("A375").Value = "PAYE5E80BEA"
$ws.Range("A376").Value = "PAYA035742D"
$ws.Range("A377").Value = "PAYCCA63A97"
$ws.Range("A378").Value = "PAY0F9F913E"
$ws.Range("A379").Value = "PAYF77BFC3D"
$ws.Range("A380").Value = "PAY1BA12561"
$ws.Range("A381").Value = "PAYA21A273F"
$ws.Range("A382").Value = "PAY2DAC2FD5"
$ws.Range("A383").Value = "PAY112B28F4"
$ws.Range("A384").Value = "PAY4023A2CA"
$ws.Range("A385").Value = "PAY55C85DEF"
$ws.Range("A386").Value = "PAY5F07155B"
$ws.Range("A387").Value = "PAY3FFB8B30"
$ws.Range("A388").Value = "PAY194F3A2F"
$ws.Range("A389").Value = "PAY3D0A2B01"
$ws.Range("A390").Value = "PAY217D9A63"
$ws.Range("A391").Value = "PAY5812F30D"
$ws.Range("A392").Value = "PAY7A54C24A"
$ws.Range("A393").Value = "PAYC089177A"
$ws.Range("A394").Value = "PAY7E58445E"
$ws.Range("A395").Value = "PAYB2635D51"
$ws.Range("A396").Value = "PAY6AC060C6"
$ws.Range("A397").Value = "PAY11324923"
$ws.Range("A398").Value = "PAY7EB3D9F8"
$ws.Range("A399").Value = "PAYA6B42B4B"
$ws.Range("A400").Value = "PAY259784C9"
$ws.Range("A401").Value = "PAYCAE43906"
$ws.Range("A402").Value = "PAY0C5E7177"
$ws.Range("A403").Value = "PAYEF6F74A2"
$ws.Range("A404").Value = "PAY345CE79B"
$ws.Range("A405").Value = "PAY56325639"
$ws.Range("A406").Value = "PAYE701F47C"
$ws.Range("A407").Value = "PAYBDC32C20"
$ws.Range("A408").Value = "PAY21D8882B"
$ws.Range("A409").Value = "PAYCB6C9C38"
$ws.Range("A410").Value = "PAY3E2D9BCA"
$ws.Range("A411").Value = "PAY31AB27CF"
$ws.Range("A412").Value = "PAYCA98FB27"
$ws.Range("A413").Value = "PAY2C766539"
$ws.Range("A414").Value = "PAY91180446"
$ws.Range("A415").Value = "PAYD48E4AD4"
$ws.Range("A416").Value = "PAYE731A849"
$ws.Range("A417").Value = "PAY8794B867"
$ws.Range("A418").Value = "PAY6D0317FD"
$ws.Range("A419").Value = "PAYA6B69F9D"
$ws.Range("A420").Value = "PAY3214822F"
$ws.Range("A421").Value = "PAYBA949F5C"
$ws.Range("A422").Value = "PAY76BEE84B"
$ws.Range("A423").Value = "PAYB0D802F0"
$ws.Range("A424").Value = "PAYFFBB83C4"
$ws.Range("A425").Value = "PAY5DE9C4E1"
$ws.Range("A426").Value = "PAY77E3BD52"
$ws.Range("A427").Value = "PAY6C233289"
$ws.Range("A428").Value = "PAYD7343F6D"
$ws.Range("A429").Value = "PAY71883CCA"
$ws.Range("A430").Value = "PAY02C71E79"
$ws.Range("A431").Value = "PAY1DE94E77"
$ws.Range("A432").Value = "PAY6A8AA51C"
$ws.Range("A433").Value = "PAYCD286434"
$ws.Range("A434").Value = "PAY63B0A77E"
$ws.Range("A435").Value = "PAYB823DD23"
$ws.Range("A436").Value = "PAY8C19E5E2"
$ws.Range("A437").Value = "PAYCCCF8172"
$ws.Range("A438").Value = "PAY8AB0731B"
$ws.Range("A439").Value = "PAY2F13D20A"
$ws.Range("A440").Value = "PAY82950469"
$ws.Range("A441").Value = "PAY280BC052"
$ws.Range("A442").Value = "PAYC0634CE4"
$ws.Range("A443").Value = "PAY3B2317AC"
$ws.Range("A444").Value = "PAY0F3DA786"
$ws.Range("A445").Value = "PAY9ECD1EC5"
$ws.Range("A446").Value = "PAY3D208CA2"
$ws.Range("A447").Value = "PAYC2C60572"
$ws.Range("A448").Value = "PAYC75C623D"
$ws.Range("A449").Value = "PAY8EFAA1F2"
$ws.Range("A450").Value = "PAYB7D162A6"
$ws.Range("A451").Value = "PAYB1920ABB"
$ws.Range("A452").Value = "PAY876BDD7A"
$ws.Range("A453").Value = "PAY0F82F605"
$ws.Range("A454").Value = "PAYD815ADA9"
$ws.Range("A455").Value = "PAYCB7F162D"
$ws.Range("A456").Value = "PAY61842528"
$ws.Range("A457").Value = "PAY6F5B5987"
$ws.Range("A458").Value = "PAY77714833"
$ws.Range("A459").Value = "PAYF8ABFCE0"
$ws.Range("A460").Value = "PAYC39725DC"
$ws.Range("A461").Value = "PAY3FD80A9F"
$ws.Range("A462").Value = "PAY8C54E490"
$ws.Range("A463").Value = "PAYAFE9B612"
$ws.Range("A464").Value = "PAY395F7656"
$ws.Range("A465").Value = "PAY32707304"
$ws.Range("A466").Value = "PAYFC0E2891"
$ws.Range("A467").Value = "PAY4F885397"
$ws.Range("A468").Value = "PAY2049EADE"
$ws.Range("A469").Value = "PAYE46610E4"
$ws.Range("A470").Value = "PAY1A709E55"
$ws.Range("A471").Value = "PAYA14F24EA"
$ws.Range("A472").Value = "PAY1D96F91C"
$ws.Range("A473").Value = "PAY800A119E"
$ws.Range("A474").Value = "PAYD1640D0F"
$ws.Range("A475").Value = "PAY4DAC1733"
$ws.Range("A476").Value = "PAYFB38A771"
$ws.Range("A477").Value = "PAY6FDE1C5A"
$ws.Range("A478").Value = "PAY4808353A"
$ws.Range("A479").Value = "PAY87A592FF"
$ws.Range("A480").Value = "PAY7A02DED8"
$ws.Range("A481").Value = "PAY5D653CF9"
$ws.Range("A482").Value = "PAYC2CDD9C2"
$ws.Range("A483").Value = "PAYB7653FE3"
$ws.Range("A484").Value = "PAY50DED40D"
$ws.Range("A485").Value = "PAYD90082F7"
$ws.Range("A486").Value = "PAYB3706E98"
$ws.Range("A487").Value = "PAY25F72069"
$ws.Range("A488").Value = "PAY861979AB"
$ws.Range("A489").Value = "PAYEAF5E282"
$ws.Range("A490").Value = "PAY6F55A867"
$ws.Range("A491").Value = "PAYE7C1508D"
